$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Unprotect()

# Update the confidential footer date from 2021-03-30 to 2021-03-31
$ws.Range("A10").Value = "***CONFIDENTIAL***: For one-on-one client use only. Not approved for distribution.`nModel holdings provided as of 2021-03-31 for illustrative purposes only and are subject to change."

# Update the Weight (D) and Percent Change (E) values for rows 2-7
$ws.Range("D2").Value = 0.4886695214274901
$ws.Range("E2").Value = 0.005122143420015934

$ws.Range("D3").Value = 0.3340564183666554
$ws.Range("E3").Value = -0.003352065463866793

$ws.Range("D4").Value = 0.09381315778306562
$ws.Range("E4").Value = -0.009451795841209809

$ws.Range("D5").Value = 0.05495024013812296
$ws.Range("E5").Value = 0.004840940525587811

$ws.Range("D6").Value = 0.0285106622846658
$ws.Range("E6").Value = 0.01632789070309881

$ws.Range("E7").Value = 0.001228083398199908

# Restore sheet protection (it was removed above to allow edits on the protected sheet)
$ws.Protect()
